# Update the crypto price (D) and 1h volume change (E) columns
# with freshly scraped figures (GitHub Actions daily refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) ---
# A leading apostrophe forces Excel to store the value as text instead
# of re-interpreting number-like strings (e.g. "58.14") as numeric
# values, matching how these price cells were originally stored (text).
$ws.Range("D2").Value = '38.221.82'
$ws.Range("D3").Value = '2.071.64'
$ws.Range("D6").Value = '''0.617'
$ws.Range("D7").Value = '''58.14'
$ws.Range("D9").Value = '''0.387'
$ws.Range("D12").Value = '2.379.02'
$ws.Range("D13").Value = '''14.63'
$ws.Range("D14").Value = '''20.79'
$ws.Range("D17").Value = '2.067.28'
$ws.Range("D18").Value = '38.162.22'
$ws.Range("D19").Value = '''6.25'
$ws.Range("D20").Value = '''70.04'
$ws.Range("D21").Value = '0.0₃0832'
$ws.Range("D22").Value = '''225.27'
$ws.Range("D26").Value = '''9.31'
$ws.Range("D27").Value = '''165.94'
$ws.Range("D29").Value = '''19.13'
$ws.Range("D30").Value = '''1.38'
$ws.Range("D32").Value = '''4.57'
$ws.Range("D33").Value = '''4.60'
$ws.Range("D34").Value = '''0.0616'
$ws.Range("D35").Value = '''1.99'
$ws.Range("D36").Value = '''2.40'
$ws.Range("D37").Value = '''6.08'
$ws.Range("D40").Value = '''98.55'
$ws.Range("D41").Value = '''0.0219'
$ws.Range("D42").Value = '1.485.08'
$ws.Range("D43").Value = '''0.0951'
$ws.Range("D44").Value = '''16.83'
$ws.Range("D45").Value = '''2.86'
$ws.Range("D47").Value = '''4.06'
$ws.Range("D50").Value = '''7.12'
$ws.Range("D51").Value = '2.264.44'

# --- 1h volume change (column E) ---
$ws.Range("E2").Value = '  +3.07%  '
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("E6").Value = '  +1.73%  '
$ws.Range("E7").Value = '  +6.16%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("E10").Value = '  +2.71%  '
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("E12").Value = '  +2.73%  '
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("E15").Value = '  +1.59%  '
$ws.Range("E16").Value = '  +2.69%  '
$ws.Range("E17").Value = '  +2.80%  '
$ws.Range("E18").Value = '  +3.14%  '
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("E20").Value = '  +1.74%  '
$ws.Range("E21").Value = '  +1.42%  '
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("E25").Value = '  +2.98%  '
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("E28").Value = '  +7.60%  '
$ws.Range("E29").Value = '  +2.28%  '
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("E33").Value = '  +4.36%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  +7.75%  '
$ws.Range("E36").Value = '  +2.19%  '
$ws.Range("E37").Value = '  +12.48%  '
$ws.Range("E38").Value = '  +5.10%  '
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("E40").Value = '  +3.54%  '
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E43").Value = '  +2.97%  '
$ws.Range("E44").Value = '  +1.75%  '
$ws.Range("E45").Value = '  +3.74%  '
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("E47").Value = '  +15.82%  '
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("E51").Value = '  +2.62%  '
